# Update cryptos list values as per the data refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.060.27"

$ws.Range("D3").Value = "3.193.30"
$ws.Range("E3").Value = "  +1.83%  "

$ws.Range("D5").Value = "'537.82"
$ws.Range("E5").Value = "  +0.71%  "

$ws.Range("D6").Value = "'144.82"
$ws.Range("E6").Value = "  +4.07%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "'0.525"
$ws.Range("E8").Value = "  +3.92%  "

$ws.Range("E9").Value = "  -0.15%  "

$ws.Range("E10").Value = "  +4.39%  "

$ws.Range("D11").Value = "'0.430"
$ws.Range("E11").Value = "  +2.90%  "

$ws.Range("D12").Value = "3.745.86"
$ws.Range("E12").Value = "  +1.87%  "

$ws.Range("E13").Value = "  -1.14%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'26.16"
$ws.Range("E14").Value = "  +1.93%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0000175"
$ws.Range("E15").Value = "  +2.23%  "

$ws.Range("D16").Value = "60.105.36"

$ws.Range("D17").Value = "3.198.27"
$ws.Range("E17").Value = "  +2.02%  "

$ws.Range("D18").Value = "'6.21"
$ws.Range("E18").Value = "  -0.54%  "

$ws.Range("D19").Value = "'13.10"
$ws.Range("E19").Value = "  +1.26%  "

$ws.Range("D20").Value = "'8.39"
$ws.Range("E20").Value = "  +2.64%  "

$ws.Range("D21").Value = "'383.78"
$ws.Range("E21").Value = "  +2.28%  "

$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  -0.24%  "

$ws.Range("E23").Value = "  +3.05%  "

$ws.Range("D24").Value = "'70.29"
$ws.Range("E24").Value = "  +0.07%  "

$ws.Range("D25").Value = "'0.171"
$ws.Range("E25").Value = "  +2.12%  "

$ws.Range("D26").Value = "'8.76"
$ws.Range("E26").Value = "  +11.74%  "

$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").Value = "0.0₃0899"
$ws.Range("E28").Value = "  +1.56%  "

$ws.Range("E29").Value = "  +1.31%  "

$ws.Range("D30").Value = "'22.36"
$ws.Range("E30").Value = "  +2.73%  "

$ws.Range("D31").Value = "'6.16"
$ws.Range("E31").Value = "  -0.87%  "

$ws.Range("E32").Value = "  +3.89%  "

$ws.Range("E33").Value = "  +3.10%  "

$ws.Range("D34").Value = "'6.53"
$ws.Range("E34").Value = "  +4.64%  "

$ws.Range("D35").Value = "'156.43"
$ws.Range("E35").Value = "  -2.93%  "

$ws.Range("E36").Value = "  +0.29%  "

$ws.Range("D37").Value = "2.763.56"
$ws.Range("E37").Value = "  +5.31%  "

$ws.Range("D38").Value = "'25.58"
$ws.Range("E38").Value = "  +0.22%  "

$ws.Range("D39").Value = "'0.0711"
$ws.Range("E39").Value = "  +5.58%  "

$ws.Range("E40").Value = "  +0.14%  "

$ws.Range("D41").Value = "'4.27"
$ws.Range("E41").Value = "  +1.18%  "

$ws.Range("D42").Value = "'39.78"
$ws.Range("E42").Value = "  +1.94%  "

$ws.Range("D43").Value = "'0.728"
$ws.Range("E43").Value = "  +4.01%  "

$ws.Range("E44").Value = "  +5.57%  "

$ws.Range("D45").Value = "3.235.42"
$ws.Range("E45").Value = "  +1.77%  "

$ws.Range("E46").Value = "  +2.54%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.101"
$ws.Range("E47").Value = "  +1.90%  "

$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "'6.19"
$ws.Range("E48").Value = "  -0.57%  "

$ws.Range("D49").Value = "'0.797"
$ws.Range("E49").Value = "  +6.05%  "

$ws.Range("D50").Value = "'20.47"
$ws.Range("E50").Value = "  +1.04%  "
